# "Bugfixed the naive forecaster component module"
#
# The `date` column (A2:A22) currently holds raw date serials (Dec-31 of
# each year, 2004-2024) formatted with a custom "YYYY-MM-DD HH:MM:SS"
# number format. The forecaster should instead treat these as quarter
# labels, so replace each date with its "<year>Q4" text label and give it
# the same (plain, non-date) style already used by the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$firstDataRow = 2
$lastDataRow = 22

# Row N (2..22) holds Dec-31 of year (2002 + N) -> label "<year>Q4".
for ($row = $firstDataRow; $row -le $lastDataRow; $row++) {
    $year = 2002 + $row
    $ws.Cells.Item($row, 1).Value = "$($year)Q4"
}

# Re-style column A's data cells to match the header's plain text style
# (centered/bold, bordered, General format) instead of the old date-time
# number format, by copying the header cell's format down.
$ws.Cells.Item(1, 1).Copy()
$ws.Range($ws.Cells.Item($firstDataRow, 1), $ws.Cells.Item($lastDataRow, 1)).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The custom date-time number format is no longer used anywhere in the
# workbook now that column A holds text labels - drop it.
$wb.DeleteNumberFormat("YYYY-MM-DD HH:MM:SS")
